# Apply the commit's data changes to the workbook.
$wb = $excel.ActiveWorkbook

# --- 1) available_players sheet ---
# The roster had a stray duplicate "Holloran" entry at the bottom (row 67) and
# was missing a newly-evaluated player, "Acconcia, Sienna", who alphabetically
# belongs right before "Assan, Jonnelle" (row 27). Shift the Player column
# (C27:C66) down one row into C28:C67 -- this both drops the stray "Holloran"
# duplicate off the end and makes room at the top -- then write in the new
# player's name. Eval Number / Grade / MW columns (A, B, D) are untouched.
$wsPlayers = $wb.Worksheets.Item("available_players")
$wsPlayers.Range("C27:C66").Copy($wsPlayers.Range("C28:C67"))
$wsPlayers.Range("C27").Value2 = "Acconcia, Sienna"

# Leave the selection where the edit landed, on the newly-inserted player.
$wsPlayers.Activate() | Out-Null
$wsPlayers.Range("C28").Select() | Out-Null

# --- 2) draft_order sheet ---
# Populate the previously-empty draft_order sheet with the snake-draft order.
$wsOrder = $wb.Worksheets.Item("draft_order")
$wsOrder.Range("A1").Value2 = "Order"
$wsOrder.Range("A2").Value2 = "shanks"
$wsOrder.Range("A3").Value2 = "hirsch"
$wsOrder.Range("A4").Value2 = "riley"
$wsOrder.Range("A5").Value2 = "gianarikas"
$wsOrder.Range("A6").Value2 = "barret"
$wsOrder.Range("A7").Value2 = "hurley"
$wsOrder.Range("A8").Value2 = "baker"
$wsOrder.Range("A9").Value2 = "oriely"

# Leave the selection on the next empty row, like after typing a list.
$wsOrder.Activate() | Out-Null
$wsOrder.Range("A10").Select() | Out-Null
